# Update cryptocurrency price/volume data (symbol list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> D (Price) and E (Volume 1h) new values.
# Some rows only need the E column updated (D unchanged).
$updates = @(
    @{ Row = 2;  D = "287.27";    E = "1.65%" }
    @{ Row = 3;  D = "29.15";     E = "2.43%" }
    @{ Row = 4;  D = "5.214";     E = "3.36%" }
    @{ Row = 5;  D = "0.06978";   E = "6.76%" }
    @{ Row = 6;  D = "7.412";     E = "1.97%" }
    @{ Row = 7;  D = "3.560";     E = "5.75%" }
    @{ Row = 8;  D = "1.402";     E = "2.31%" }
    @{ Row = 9;  D = "0.8952";    E = "-3.62%" }
    @{ Row = 10; D = "0.1608";    E = "3.54%" }
    @{ Row = 11; D = "0.07660";   E = "27.58%" }
    @{ Row = 12; D = "0.07722";   E = "1.84%" }
    @{ Row = 13; D = "0.02923";   E = "1.02%" }
    @{ Row = 14; D = "0.09004";   E = "0.44%" }
    @{ Row = 15; D = "0.001586";  E = "0.52%" }
    @{ Row = 16; D = "0.0006480"; E = "1.89%" }
    @{ Row = 17; D = "0.006437";  E = "6.19%" }
    @{ Row = 18; D = "3.460";     E = "0.41%" }
    @{ Row = 19; D = "2.230";     E = "-0.24%" }
    @{ Row = 20; D = $null;       E = "0.88%" }
    @{ Row = 21; D = "0.1330";    E = "4.50%" }
    @{ Row = 22; D = "4.008";     E = "-1.97%" }
    @{ Row = 23; D = "0.1551";    E = "1.92%" }
    @{ Row = 24; D = "0.04522";   E = "1.31%" }
    @{ Row = 25; D = "0.001207";  E = "2.46%" }
    @{ Row = 26; D = "0.004240";  E = "-3.84%" }
    @{ Row = 27; D = "0.0001166"; E = "-6.42%" }
    @{ Row = 28; D = "0.0001621"; E = "0.58%" }
    @{ Row = 40; D = "0.04339";   E = "4.41%" }
    @{ Row = 41; D = "0.006898";  E = "4.41%" }
    @{ Row = 42; D = $null;       E = "1.61%" }
    @{ Row = 43; D = "0.002072";  E = "2.95%" }
    @{ Row = 44; D = "0.01175";   E = "-2.26%" }
    @{ Row = 45; D = "0.00005825"; E = "5.71%" }
    @{ Row = 47; D = "0.01304";   E = "0.51%" }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cellD = $ws.Cells.Item($r, 4)
        # Force text storage (these look numeric, so without this Excel
        # would coerce them to the Number type), then restore the
        # original "General" display format used by the rest of the sheet.
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.D
        $cellD.NumberFormat = "General"
    }
    $cellE = $ws.Cells.Item($r, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $u.E
    $cellE.NumberFormat = "General"
}
